$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.594.21"
$ws.Range("E2").Value = "  +2.33%  "
$ws.Range("D3").Value = "2.408.45"
$ws.Range("E3").Value = "  +2.81%  "
$ws.Range("D5").Value = "'552.71"
$ws.Range("E5").Value = "  +1.98%  "
$ws.Range("D6").Value = "'136.89"
$ws.Range("E6").Value = "  +1.52%  "
$ws.Range("D8").Value = "'0.569"
$ws.Range("E8").Value = "  +1.35%  "
$ws.Range("E9").Value = "  +5.21%  "
$ws.Range("D10").Value = "'5.79"
$ws.Range("E10").Value = "  +3.29%  "
$ws.Range("D11").Value = "'0.362"
$ws.Range("E11").Value = "  +1.79%  "
$ws.Range("E12").Value = "  -2.13%  "
$ws.Range("D13").Value = "'24.63"
$ws.Range("E13").Value = "  +3.37%  "
$ws.Range("D14").Value = "2.839.06"
$ws.Range("E14").Value = "  +2.98%  "
$ws.Range("D15").Value = "59.519.83"
$ws.Range("E15").Value = "  +2.31%  "
$ws.Range("D16").Value = "'0.0000139"
$ws.Range("E16").Value = "  +4.33%  "
$ws.Range("D17").Value = "2.404.67"
$ws.Range("E17").Value = "  +3.89%  "
$ws.Range("D18").Value = "'11.33"
$ws.Range("E18").Value = "  +5.85%  "
$ws.Range("D19").Value = "'4.45"
$ws.Range("E19").Value = "  +4.65%  "
$ws.Range("D20").Value = "'335.53"
$ws.Range("E20").Value = "  +0.70%  "
$ws.Range("D21").Value = "'6.98"
$ws.Range("E21").Value = "  +4.54%  "
$ws.Range("D22").Value = "'1.00"
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("D23").Value = "'64.72"
$ws.Range("E23").Value = "  +2.89%  "
$ws.Range("E24").Value = "  +0.81%  "
$ws.Range("D25").Value = "'8.47"
$ws.Range("E25").Value = "  -0.83%  "
$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("D27").Value = "'1.38"
$ws.Range("E27").Value = "  -2.20%  "
$ws.Range("D28").Value = "0.0₃0778"
$ws.Range("E28").Value = "  +5.79%  "
$ws.Range("D29").Value = "'1.80"
$ws.Range("E29").Value = "  +1.98%  "
$ws.Range("D30").Value = "'170.69"
$ws.Range("E30").Value = "  -0.08%  "
$ws.Range("D31").Value = "'6.27"
$ws.Range("E31").Value = "  +2.62%  "
$ws.Range("D32").Value = "'18.70"
$ws.Range("E32").Value = "  +1.42%  "
$ws.Range("D33").Value = "'1.02"
$ws.Range("E33").Value = "  -0.72%  "
$ws.Range("E34").Value = "  -0.01%  "
$ws.Range("D35").Value = "'4.31"
$ws.Range("E35").Value = "  +1.24%  "
$ws.Range("D36").Value = "'1.31"
$ws.Range("E36").Value = "  +5.24%  "
$ws.Range("E37").Value = "  +0.09%  "
$ws.Range("E38").Value = "  -0.47%  "
$ws.Range("D39").Value = "'40.10"
$ws.Range("E39").Value = "  +2.62%  "
$ws.Range("D40").Value = "'0.419"
$ws.Range("E40").Value = "  +11.39%  "
$ws.Range("D41").Value = "'303.78"
$ws.Range("E41").Value = "  +6.17%  "
$ws.Range("D42").Value = "'3.76"
$ws.Range("E42").Value = "  +3.06%  "
$ws.Range("D43").Value = "'142.36"
$ws.Range("E43").Value = "  -0.50%  "
$ws.Range("B44").Value = "Stellar"
$ws.Range("C44").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D44").Value = "'0.0963"
$ws.Range("E44").Value = "  +2.54%  "
$ws.Range("B45").Value = "Hedera"
$ws.Range("C45").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D45").Value = "'0.0526"
$ws.Range("E45").Value = "  +4.43%  "
$ws.Range("B46").Value = "Polygon"
$ws.Range("C46").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D46").Value = "'0.412"
$ws.Range("E46").Value = "  +5.93%  "
$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").Value = "'0.571"
$ws.Range("E47").Value = "  +1.29%  "
$ws.Range("E48").Value = "  -0.72%  "
$ws.Range("D49").Value = "'0.0226"
$ws.Range("E49").Value = "  +3.63%  "
$ws.Range("D50").Value = "'11.03"
$ws.Range("E50").Value = "  -0.40%  "
$ws.Range("E51").Value = "  +4.39%  "
